# fixed beril-model DH Playground link
# Add two new worksheets ("Splitting" and "Pooling") after the existing
# "Process" sheet, each with the standard Process-style header row:
# inputs | outputs | id | name | description

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSplitting = $wb.Worksheets.Add($null, $lastSheet)
$wsSplitting.Name = "Splitting"
$wsSplitting.Range("A1").Value = "inputs"
$wsSplitting.Range("B1").Value = "outputs"
$wsSplitting.Range("C1").Value = "id"
$wsSplitting.Range("D1").Value = "name"
$wsSplitting.Range("E1").Value = "description"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPooling = $wb.Worksheets.Add($null, $lastSheet)
$wsPooling.Name = "Pooling"
$wsPooling.Range("A1").Value = "inputs"
$wsPooling.Range("B1").Value = "outputs"
$wsPooling.Range("C1").Value = "id"
$wsPooling.Range("D1").Value = "name"
$wsPooling.Range("E1").Value = "description"
